# Auto-generated script applying the cryptos.xlsx data refresh diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "96.678.45"
$ws.Range("E2").Value = "  -2.01%  "
$ws.Range("D3").Value = "3.299.97"
$ws.Range("E3").Value = "  -5.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.34"
$ws.Range("E5").Value = "  -6.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "644.78"
$ws.Range("E6").Value = "  -4.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.34"
$ws.Range("E7").Value = "  -13.74%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.409"
$ws.Range("E8").Value = "  -11.82%  "
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.974"
$ws.Range("E10").Value = "  -11.36%  "
$ws.Range("D11").Value = "3.299.16"
$ws.Range("E11").Value = "  -5.58%  "
$ws.Range("E12").Value = "  -8.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "39.42"
$ws.Range("E13").Value = "  -7.64%  "
$ws.Range("D14").Value = "96.489.03"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.95"
$ws.Range("E15").Value = "  -4.09%  "
$ws.Range("E16").Value = "  -9.88%  "
$ws.Range("D17").Value = "3.917.05"
$ws.Range("E17").Value = "  -5.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.55"
$ws.Range("E18").Value = "  +4.65%  "
$ws.Range("D19").Value = "3.300.33"
$ws.Range("E19").Value = "  -5.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.48"
$ws.Range("E20").Value = "  -6.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.489"
$ws.Range("E21").Value = "  +8.41%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.30"
$ws.Range("E22").Value = "  -4.96%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "486.66"
$ws.Range("E23").Value = "  -9.54%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.22"
$ws.Range("E24").Value = "  -11.63%  "
$ws.Range("E25").Value = "  -12.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.27"
$ws.Range("E26").Value = "  -2.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "93.64"
$ws.Range("E27").Value = "  -8.80%  "
$ws.Range("E28").Value = "  -9.16%  "
$ws.Range("D29").Value = "3.476.23"
$ws.Range("E29").Value = "  -5.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -6.36%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.65"
$ws.Range("E32").Value = "  -6.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.185"
$ws.Range("E33").Value = "  -7.80%  "
$ws.Range("E34").Value = "  +8.20%  "
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.536"
$ws.Range("E36").Value = "  -8.81%  "
$ws.Range("E37").Value = "  -10.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.44"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.37"
$ws.Range("E39").Value = "  -8.18%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("E41").Value = "  -8.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "493.87"
$ws.Range("E42").Value = "  -9.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.50"
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("E44").Value = "  -2.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.810"
$ws.Range("E45").Value = "  -7.35%  "
$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.32"
$ws.Range("E46").Value = "  +0.60%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0400"
$ws.Range("E47").Value = "  -8.58%  "
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.34"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.60"
$ws.Range("E49").Value = "  +0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.17"
$ws.Range("E50").Value = "  +1.47%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.07"
$ws.Range("E51").Value = "  -12.99%  "
